# Rename the sheet from its literal month code ("11005") to the generic
# "YYYMM" placeholder, matching the already-generic title text in A1
# ("YYY.MM  轉催收明細總表"). Renaming the sheet automatically repoints the
# '_xlnm._FilterDatabase' defined name; the '_xlnm.Print_Area' defined name
# is refreshed explicitly via PageSetup.PrintArea so it also points at the
# renamed sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "YYYMM"
$ws.PageSetup.PrintArea = "`$A`$1:`$M`$10"
